# Updated cryptos list with GitHub Actions (refreshed Price / Volume(1h) data).
# Rows 42/43 also swap the Algorand <-> RenderToken entries to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.667.73"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.922.83"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "'339.31"
$ws.Range("E5").Value = "  +4.30%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "'0.4812"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.08092"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'1.001"
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").Value = "'23.39"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "1.987.11"
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "'5.987"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "'7.189"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "'90.04"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "'0.06850"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "'1.014"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.00001028"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'17.53"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "'1.011"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "29.677.92"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'5.566"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "'11.78"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'2.162"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "2.188.44"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'6.592"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "'157.26"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").Value = "'19.87"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "'2.065"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'120.41"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "'0.09589"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'5.525"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").Value = "'1.399"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "'3.550"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'0.06543"
$ws.Range("E36").Value = "  +6.98%  "
$ws.Range("D37").Value = "'0.02261"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "'1.200"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").Value = "'0.5896"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").Value = "'7.864"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.490"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1832"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "'1.243"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").Value = "'0.07464"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").Value = "'0.5509"
$ws.Range("D48").Value = "'1.964"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'116.43"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "'2.407"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").Value = "'71.82"
$ws.Range("E51").Value = "  -0.82%  "
